$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 7.2
$ws.Range("G2").Value = 7.4
$ws.Range("H2").Value = 1.48
$ws.Range("I2").Value = 1.51
$ws.Range("J2").Value = 4.9
$ws.Range("K2").Value = 5.2
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 2.92
$ws.Range("W2").Value = 1.15
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 1000
$ws.Range("Z2").Value = 1000
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 1000
$ws.Range("AC2").Value = 1000
$ws.Range("AD2").Value = 1000
$ws.Range("AE2").Value = 10.5
$ws.Range("AF2").Value = 1000
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 1000
$ws.Range("AI2").Value = 7
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 1000
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 8.6
$ws.Range("AN2").Value = 980
$ws.Range("AO2").Value = 2.68

# Row 3
$ws.Range("F3").Value = 3.1
$ws.Range("G3").Value = 3.2
$ws.Range("H3").Value = 2.84
$ws.Range("I3").Value = 2.88
$ws.Range("J3").Value = 2.96
$ws.Range("K3").Value = 3.05
$ws.Range("L3").Value = 3.7
$ws.Range("M3").Value = 1.17
$ws.Range("N3").Value = 2.28
$ws.Range("O3").Value = 1.76
$ws.Range("P3").Value = 1.4
$ws.Range("Q3").Value = 3.4
$ws.Range("R3").Value = 1.13
$ws.Range("S3").Value = 8
$ws.Range("T3").Value = 2.34
$ws.Range("U3").Value = 1.64
$ws.Range("V3").Value = 1.53
$ws.Range("W3").Value = 1.45
$ws.Range("X3").Value = 7.8
$ws.Range("Y3").Value = 8.199999999999999
$ws.Range("Z3").Value = 18
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 8.800000000000001
$ws.Range("AC3").Value = 7.2
$ws.Range("AD3").Value = 17.5
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 18
$ws.Range("AG3").Value = 18.5
$ws.Range("AH3").Value = 990
$ws.Range("AI3").Value = 1000
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AL3").Value = 1000
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 1000

# Row 4
$ws.Range("F4").Value = 3.4
$ws.Range("G4").Value = 3.5
$ws.Range("H4").Value = 2.38
$ws.Range("I4").Value = 2.44
$ws.Range("J4").Value = 3.35
$ws.Range("K4").Value = 3.45
$ws.Range("L4").Value = 3
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 3
$ws.Range("O4").Value = 1.48
$ws.Range("P4").Value = 1.6
$ws.Range("Q4").Value = 2.56
$ws.Range("R4").Value = 1.21
$ws.Range("S4").Value = 5.5
$ws.Range("T4").Value = 2.16
$ws.Range("U4").Value = 1.79
$ws.Range("V4").Value = 1.7
$ws.Range("W4").Value = 1.4
$ws.Range("X4").Value = 11.5
$ws.Range("Y4").Value = 7.6
$ws.Range("Z4").Value = 12
$ws.Range("AA4").Value = 28
$ws.Range("AB4").Value = 11.5
$ws.Range("AC4").Value = 7.4
$ws.Range("AD4").Value = 12
$ws.Range("AE4").Value = 30
$ws.Range("AF4").Value = 25
$ws.Range("AG4").Value = 17.5
$ws.Range("AH4").Value = 990
$ws.Range("AI4").Value = 80
$ws.Range("AJ4").Value = 90
$ws.Range("AK4").Value = 980
$ws.Range("AL4").Value = 95
$ws.Range("AM4").Value = 220
$ws.Range("AN4").Value = 90
$ws.Range("AO4").Value = 980

# Row 5
$ws.Range("F5").Value = 9
$ws.Range("G5").Value = 9.199999999999999
$ws.Range("H5").Value = 1.41
$ws.Range("I5").Value = 1.42
$ws.Range("J5").Value = 5.3
$ws.Range("K5").Value = 5.6
$ws.Range("L5").Value = 1.32
$ws.Range("N5").Value = 5.6
$ws.Range("O5").Value = 1.21
$ws.Range("P5").Value = 2.52
$ws.Range("Q5").Value = 1.62
$ws.Range("R5").Value = 1.58
$ws.Range("S5").Value = 2.66
$ws.Range("T5").Value = 1.81
$ws.Range("U5").Value = 2.14
$ws.Range("V5").Value = 3.1
$ws.Range("W5").Value = 1.13
$ws.Range("X5").Value = 25
$ws.Range("Y5").Value = 10.5
$ws.Range("Z5").Value = 10
$ws.Range("AA5").Value = 13
$ws.Range("AB5").Value = 34
$ws.Range("AE5").Value = 14.5
$ws.Range("AF5").Value = 80
$ws.Range("AG5").Value = 990
$ws.Range("AH5").Value = 22
$ws.Range("AI5").Value = 29
$ws.Range("AJ5").Value = 250
$ws.Range("AK5").Value = 90
$ws.Range("AL5").Value = 110
$ws.Range("AM5").Value = 110
$ws.Range("AN5").Value = 130
$ws.Range("AO5").Value = 6

# Row 6
$ws.Range("G6").Value = 3.3
$ws.Range("H6").Value = 2.44
$ws.Range("K6").Value = 4.4
$ws.Range("T6").Value = 1.05
$ws.Range("AB6").Value = 46
$ws.Range("AC6").Value = 19
$ws.Range("AD6").Value = 32
$ws.Range("AG6").Value = 40
